# GBDS JANUARY FILES 2026 | fliqlo@GBDS
# Fill in the daily figures for the "01,08" sheet (BEG INVTY / TOTAL SALES /
# SALES TO TRADE breakdown columns). Dependent formula cells (I, J, R, S and
# the grand-total row 52) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01,08")
$ws.Activate()

$values = @{
    "C8"  = 89
    "C9"  = 143
    "C11" = 18;  "G11" = 6;  "L11" = 1;  "P11" = 5
    "C12" = 33;  "G12" = 5;  "P12" = 5
    "C13" = 3957; "D13" = 9; "G13" = 42; "H13" = 20
    "L13" = 13;  "M13" = 3;  "N13" = 15; "O13" = 11; "P13" = 14; "Q13" = 6
    "C14" = 5
    "C15" = 1
    "C16" = 1
    "C18" = 6
    "C20" = 4
    "C21" = 7
    "C22" = 6
    "C23" = 4
    "C24" = 55
    "C25" = 8
    "C26" = 66; "D26" = 11; "H26" = 12; "M26" = 12
    "C27" = 21; "G27" = 1;  "L27" = 1
    "C28" = 83; "D28" = 12; "G28" = 3;  "H28" = 12; "M28" = 12; "P28" = 3
    "C29" = 10; "G29" = 1;  "P29" = 1
    "C30" = 102
    "C33" = 5
    "C34" = 209; "D34" = 23; "G34" = 1; "P34" = 1
    "C36" = 1
    "C37" = 6917; "D37" = 2; "G37" = 114; "L37" = 10; "N37" = 20; "P37" = 84
    "C38" = 35
    "C39" = 1098; "G39" = 4; "L39" = 2;  "P39" = 2
    "C40" = 2316; "G40" = 21; "L40" = 2; "P40" = 19
    "C42" = 34124; "D42" = 3; "G42" = 1353; "L42" = 523; "N42" = 450; "P42" = 380
    "C45" = 2
    "C46" = 2415; "D46" = 3; "G46" = 4; "L46" = 1; "N46" = 3
    "C48" = 1
    "C50" = 34
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# The monthly grand-total row has two manual adjustment offsets baked into
# the SUM formulas for the ENDING INVTY columns.
$ws.Range("I52").Formula = "=SUM(I8:I51)-2"
$ws.Range("J52").Formula = "=SUM(J8:J51)+48"

# Restore the cursor position reported in the saved view state.
$ws.Range("P14").Select()
